# Applies the "Fixed errors in NOAA site file to match Steve's comments" edit
# to the siteNameMapping worksheet (first sheet, content-wise the site info table).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Row 21 (omega / raw008): Keep flag was wrongly "F"; Steve confirmed row 8 is
#     omega so it should be kept ("T"). Also remove the stale yellow highlight that
#     had flagged this row for review, and update the justification note.
$ws1.Range("A21:D21").Interior.ColorIndex = -4142
$ws1.Range("F21").Value = "T"
$ws1.Range("G21").Value = "KV changed to KEEP=T because 8 is omega. Seems to have been a copy/paste error"

# --- Row 42 (e1 / leg_e): Latitude/Longitude were copy/pasted from the wrong block;
#     fix them to match the other "e" sites (rows 43-46).
$ws1.Range("D42").Value = 43.053666669999998
$ws1.Range("E42").Value = -86.262833330000007
$ws1.Range("G42").Value = "KV fixed coordinates to match other e sites, per Steve's comment"

# --- Row 54 (gvsubuoy): Per Steve's comment this site should be excluded (Keep = "F"),
#     and the justification note is expanded to explain why.
$ws1.Range("F54").Value = "F"
$ws1.Range("G54").Value = "see coordinates list (in Muskegon Lake) - KV changed to KEEP==F per Steve comment"
